$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 47 (shifts existing rows 47-144 down to 48-145)
$ws.Rows("47").Insert()

# Populate the new row 47 with the new data record
$ws.Range("A47").Value = 5
$ws.Range("B47").Value = "Macroferia Regional de Talca"
$ws.Range("C47").Value = "Maule"
$ws.Range("D47").Value = 44804
$ws.Range("E47").Value = 7
$ws.Range("F47").Value = "Fruta"
$ws.Range("G47").Value = 100108
$ws.Range("H47").Value = "Tropicales y subtropicales"
$ws.Range("I47").Value = 100108002
$ws.Range("J47").Value = "Mango"
$ws.Range("K47").Value = "Sin especificar"
$ws.Range("L47").Value = "Primera"
$ws.Range("M47").Value = 228
$ws.Range("N47").Value = 10000
$ws.Range("O47").Value = 10000
$ws.Range("P47").Value = 10000
$ws.Range("Q47").Value = "$/bandeja 4 kilos"
$ws.Range("R47").Value = "Brasil"
$ws.Range("S47").Value = 2500
$ws.Range("T47").Value = 4
